$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (Rent) ---
$ws.Range("C2").Value = 45837.71226172454

# --- Update existing row 3 (was "Food"/400) to the new "Mcdonalds"/50 entry ---
$ws.Range("A3").Value = "Mcdonalds"
$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 45837.12509259259

# --- Append new row 4: Snus ---
$ws.Range("A4").Value = "Snus"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 45836.71226172454

# --- Append new row 5: food ---
$ws.Range("A5").Value = "food"
$ws.Range("B5").Value = 140
$ws.Range("C5").Value = 45835.71226172454

# --- Append new row 6: Spotify ---
$ws.Range("A6").Value = "Spotify"
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 45778.12509259259

# Make the newly-added Date column cells (C4:C6) use the same date
# number format as the existing Date column cells (copy format from C2).
$ws.Range("C2").Copy()
$ws.Range("C4:C6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Refresh the sheet dimension / ignored-errors range to cover the new rows.
$ws.Range("A1:C6").Select()
